$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the style of the other header cells (B1:E1)
$ws.Range("F1").Value = "EDAM_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New EDAM_DEF column data
$ws.Range("F2").Value = "['Alphabet for a nucleotide sequence with possible ambiguity, unknown positions and non-sequence characters.']"
$ws.Range("F3").Value = "['Alphabet for a protein sequence with possible ambiguity, unknown positions and non-sequence characters.']"
$ws.Range("F4").Value = "['Mega format for (typically aligned) sequences.']"
